# Update the lattice-multiplication exercise table: every cell's
# multiplication problem is replaced with a new one (see commit
# "Update master to output generated at c8c62b6").
#
# Each table cell holds a single run containing 5 lines separated by
# <w:br/>: the "A x B" header, the 2-digit split, the "----" rule
# (constant, never changes), and two single-digit-prefixed lines. We
# rebuild each cell's Range.Text in one shot using a vertical-tab
# (chr(11)) as the line separator, which this host maps back to
# <w:br/> on write-back while preserving the run's rPr (sz=32).
#
# NOTE: Find.Execute() in this host always operates over the whole
# document story regardless of which Range it is invoked on, so it
# cannot be used to make per-cell replacements here (many lines like
# "7|    |" repeat across cells). Writing Cell.Range.Text directly is
# correctly scoped to just that cell, so that's what we use.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11

function Set-Cell($row, $col, $line1, $line2, $line4, $line5) {
    $text = $line1 + $vt + $line2 + $vt + "  ----" + $vt + $line4 + $vt + $line5
    $t.Cell($row, $col).Range.Text = $text
}

Set-Cell 1 1 "24 x 62" "  6    2" "2|    |" "4|    |"
Set-Cell 1 2 "76 x 43" "  4    3" "7|    |" "6|    |"
Set-Cell 1 3 "78 x 76" "  7    6" "7|    |" "8|    |"

Set-Cell 2 1 "44 x 45" "  4    5" "4|    |" "4|    |"
Set-Cell 2 2 "98 x 94" "  9    4" "9|    |" "8|    |"
Set-Cell 2 3 "36 x 15" "  1    5" "3|    |" "6|    |"

Set-Cell 3 1 "77 x 88" "  8    8" "7|    |" "7|    |"
Set-Cell 3 2 "56 x 45" "  4    5" "5|    |" "6|    |"
Set-Cell 3 3 "70 x 57" "  5    7" "7|    |" "0|    |"

Set-Cell 4 1 "13 x 46" "  4    6" "1|    |" "3|    |"
Set-Cell 4 2 "94 x 29" "  2    9" "9|    |" "4|    |"
Set-Cell 4 3 "37 x 41" "  4    1" "3|    |" "7|    |"

Set-Cell 5 1 "43 x 92" "  9    2" "4|    |" "3|    |"
Set-Cell 5 2 "24 x 90" "  9    0" "2|    |" "4|    |"
Set-Cell 5 3 "76 x 77" "  7    7" "7|    |" "6|    |"

Write-Host "Updated all 15 lattice multiplication cells."
